$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / URL / percentage cells - safe to assign directly
$ws.Range('D2').Value = '72.276.41'
$ws.Range('E2').Value = '  +1.89%  '
$ws.Range('D3').Value = '2.668.60'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.38%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('E9').Value = '  +5.96%  '
$ws.Range('B10').Value = 'LidoStakedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D10').Value = '2.670.36'
$ws.Range('E10').Value = '  +2.15%  '
$ws.Range('E11').Value = '  +2.18%  '
$ws.Range('E12').Value = '  +2.93%  '
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').Value = '3.159.66'
$ws.Range('E14').Value = '  +1.98%  '
$ws.Range('E15').Value = '  +3.65%  '
$ws.Range('D16').Value = '72.236.37'
$ws.Range('E16').Value = '  +1.83%  '
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = '2.672.07'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('E19').Value = '  +4.46%  '
$ws.Range('E20').Value = '  +3.37%  '
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('E22').Value = '  +1.71%  '
$ws.Range('E23').Value = '  +11.67%  '
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('B26').Value = 'NEARProtocol'
$ws.Range('C26').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('E27').Value = '  +3.88%  '
$ws.Range('D28').Value = '2.809.64'
$ws.Range('E28').Value = '  +2.77%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').Value = '0.0₃0949'
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('E31').Value = '  +2.32%  '
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('E34').Value = '  -0.21%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E36').Value = '  +0.55%  '
$ws.Range('E37').Value = '  +2.44%  '
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('E39').Value = '  +1.84%  '
$ws.Range('E40').Value = '  -6.66%  '
$ws.Range('E41').Value = '  -1.46%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('E43').Value = '  +1.05%  '
$ws.Range('E44').Value = '  -0.33%  '
$ws.Range('E45').Value = '  +1.77%  '
$ws.Range('E46').Value = '  -2.00%  '
$ws.Range('E47').Value = '  -0.53%  '
$ws.Range('E48').Value = '  +3.79%  '
$ws.Range('E49').Value = '  +3.92%  '
$ws.Range('E50').Value = '  +3.19%  '
$ws.Range('E51').Value = '  +1.60%  '

# Numeric-looking price cells that must stay TEXT (match source formatting,
# e.g. trailing zeros like "1.00"). Force text number-format before the
# assignment, then restore the original (no explicit style) afterwards so
# the cell keeps matching its un-styled siblings.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.82'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.525'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.174'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.357'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '379.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.06'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.39'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.97'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.20'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '522.44'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.59'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.39'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.07'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.335'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.33'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '153.13'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.551'
$ws.Range('D49').Style = 'Normal'
